# Timesheet update by sravani (11/03/2013)
# Mark the newly-added OFF days on the "FebruaryMarch 2013" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Rows 32-35: AP:AS all marked OFF
$ws.Range("AP32:AS32").Value = "OFF"
$ws.Range("AP33:AS33").Value = "OFF"
$ws.Range("AP34:AS34").Value = "OFF"
$ws.Range("AP35:AS35").Value = "OFF"

# Rows 36-39: AP:AR marked OFF (AS left blank)
$ws.Range("AP36:AR36").Value = "OFF"
$ws.Range("AP37:AR37").Value = "OFF"
$ws.Range("AP38:AR38").Value = "OFF"
$ws.Range("AP39:AR39").Value = "OFF"

# Reflect the user's new selection/scroll position after entering the data
$ws.Range("AS32:AS35").Select()
